# Update computed market/profit figures across the Leve profit sheets.
# This mirrors a scheduled data refresh: currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and the derived LeveProfit(NQ/HQ) columns (H:N) are
# refreshed per-row with newly pulled market data.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        [hashtable]$values,
        [array]$clear = @()
    )

    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
    foreach ($col in $clear) {
        $ws.Range("$col$row").ClearContents()
    }
}

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 17 @{
    H = 3012.2666
    I = 0
    J = 3012.2666
    K = 0
    L = 9036.799800000001
    N = -9372.799800000001
} @("M")

Set-Row $ws 96 @{
    H = 1803.3158
    I = 443.7857
    K = 1331.3571
    M = 41.64289999999983
}

Set-Row $ws 125 @{
    H = 2119.4
    I = 3199
    K = 28791
    M = -26331
}

Set-Row $ws 127 @{
    H = 877.7646999999999
    I = 876.375
    J = 900
    K = 2629.125
    L = 2700
    M = 2330.875
    N = -12620
}

Set-Row $ws 137 @{
    H = 5243.636
    I = 3575.4595
    J = 8672.666999999999
    K = 10726.3785
    L = 26018.001
    M = -8176.378499999999
    N = -31118.001
}

Set-Row $ws 141 @{
    H = 1444.931
    I = 1311.4073
    J = 3247.5
    K = 3934.2219
    L = 9742.5
    M = 1245.7781
    N = -20102.5
}

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 32 @{
    H = 19152.451
    I = 18594.418
    K = 18594.418
    M = -18307.418
}

Set-Row $ws 74 @{
    H = 325336.2
    I = 386642.47
    K = 386642.47
    M = -385768.47
}

Set-Row $ws 77 @{
    H = 325336.2
    I = 386642.47
    K = 1933212.35
    M = -1928844.35
}

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

Set-Row $ws 107 @{
    H = 1232.1428
    I = 1100.1875
    J = 1654.4
    K = 1100.1875
    L = 1654.4
    M = 819.8125
    N = -5494.4
}

Set-Row $ws 134 @{
    H = 4112.0894
    I = 3095.4482
    J = 5204.037
    K = 9286.3446
    L = 15612.111
    M = -6751.3446
    N = -20682.111
}

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

Set-Row $ws 7 @{
    H = 228.35715
    I = 215.15384
    K = 215.15384
    M = -102.15384
}

Set-Row $ws 31 @{
    H = 4891.643
    I = 2124.75
    J = 5182.8945
    K = 2124.75
    L = 5182.8945
    M = -1829.75
    N = -5772.8945
}

Set-Row $ws 34 @{
    H = 4891.643
    I = 2124.75
    J = 5182.8945
    K = 2124.75
    L = 5182.8945
    M = -1922.75
    N = -5586.8945
}

Set-Row $ws 58 @{
    H = 2703.0308
    I = 2047.2
    J = 4178.65
    K = 2047.2
    L = 4178.65
    M = -1844.2
    N = -4584.65
}

Set-Row $ws 62 @{
    H = 10772.286
    I = 4966.6665
    J = 15126.5
    K = 4966.6665
    L = 15126.5
    M = -4342.6665
    N = -16374.5
}

Set-Row $ws 65 @{
    H = 10772.286
    I = 4966.6665
    J = 15126.5
    K = 24833.3325
    L = 75632.5
    M = -21713.3325
    N = -81872.5
}

Set-Row $ws 92 @{
    H = 27557.2
    J = 27557.2
    L = 27557.2
    N = -32549.2
}

Set-Row $ws 106 @{
    H = 79999
    J = 79999
    L = 79999
    N = -82523
}

Set-Row $ws 136 @{
    H = 2703.0308
    I = 2047.2
    J = 4178.65
    K = 6141.6
    L = 12535.95
    M = -3591.6
    N = -17635.95
}

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 2 @{
    H = 2543.3572
    I = 64.5
    J = 8740.5
    K = 387
    L = 52443
    M = -274
    N = -52669
}

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

Set-Row $ws 80 @{
    H = 5657.923
    I = 4138.6
    J = 6607.5
    K = 4138.6
    L = 6607.5
    M = -3140.6
    N = -8603.5
}

Set-Row $ws 83 @{
    H = 5657.923
    I = 4138.6
    J = 6607.5
    K = 20693
    L = 33037.5
    M = -15701
    N = -43021.5
}

Set-Row $ws 132 @{
    H = 2724.935
    I = 2170.6714
    J = 8267.571
    K = 6512.014200000001
    L = 24802.713
    M = -3982.014200000001
    N = -29862.713
}

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 82 @{
    H = 1966.6666
    I = 1966.6666
    J = 0
    K = 1966.6666
    L = 0
    M = -1605.6666
} @("N")

Set-Row $ws 85 @{
    H = 1966.6666
    I = 1966.6666
    J = 0
    K = 1966.6666
    L = 0
    M = -718.6666
} @("N")

Set-Row $ws 136 @{
    H = 3374.853
    I = 2915.4138
    K = 8746.241399999999
    M = -6196.241399999999
}

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

Set-Row $ws 122 @{
    H = 1581.7106
    I = 1542.2941
    K = 4626.8823
    M = -2176.8823
}

Set-Row $ws 124 @{
    H = 109964.5
    J = 109964.5
    L = 109964.5
    N = -119784.5
}

Set-Row $ws 132 @{
    H = 2363.5303
    I = 1569.3654
    K = 4708.0962
    M = -2178.0962
}
